# Shift the timestamp column (A) back by 10800 seconds (3 hours) for every
# data row (rows 2-379); the header row (row 1) is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 379
$secondsOffset = 10800

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($null -ne $current) {
        $cell.Value = $current - $secondsOffset
    }
}
